# regen sval data to filter save games
# Updates columns B (TB), C (d2S), D (K), E (IP) and the computed G (sum)
# for rows 2-15 on the active sheet. Column A (date) and F (Win) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(3.286832544864788,  1.655778082260271,  0.1494219747398047, 10.19245300693656,  15.28448560880142)
    3  = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    4  = @(0.1190320826869504, 0.306821227259698,  0.1494219747398047, 0.4942365360607697, 1.069511820747223)
    5  = @(3.286832544864788,  10.34677158129881,  3.537761648806719,  10.19245300693656,  27.36381878190688)
    6  = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    7  = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    8  = @(3.286832544864788,  1.655778082260271,  0.1494219747398047, 0.4942365360607697, 5.586269137925634)
    9  = @(0.2917716402565462, 0.306821227259698,  3.537761648806719,  0.4942365360607697, 4.630591052383734)
    10 = @(1.455362044514542,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    11 = @(1.455362044514542,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 4.358119930609447)
    12 = @(0.1190320826869504, 0.002571899574220771, 0.7527432677738641, 0.4942365360607697, 1.368583786095805)
    13 = @(3.286832544864788,  1.655778082260271,  0.7527432677738641, 0.4942365360607697, 6.189590430959694)
    14 = @(1.455362044514542,  1.655778082260271,  3.537761648806719,  0.4942365360607697, 7.143138311642302)
    15 = @(0.6606524410359556, 1.655778082260271,  22.3905356188092,   10.19245300693656,  34.89941914904198)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]   # B - TB
    $ws.Cells.Item($row, 3).Value = $vals[1]   # C - d2S
    $ws.Cells.Item($row, 4).Value = $vals[2]   # D - K
    $ws.Cells.Item($row, 5).Value = $vals[3]   # E - IP
    $ws.Cells.Item($row, 7).Value = $vals[4]   # G - sum
}
